$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New comment for Boyd Gaming (row 2) in column T
$ws.Range("T2").Value = "ex growth, vermutlich halbes Jahr zu frühes Timing"

# Updated price for International Game Technology (row 5, IGT)
$ws.Range("I5").Value = 20.27

# Replace hardcoded PE1/PE2/PEG1/PEG2 with live formulas for row 5
$ws.Range("O5").Formula = "=I5/K5"
$ws.Range("P5").Formula = "=I5/L5"
$ws.Range("Q5").Formula = "=O5/(M5*100)"
$ws.Range("R5").Formula = "=P5/(N5*100)"
$ws.Range("O5:Q5").NumberFormat = "0.00"

# Move the active selection cursor to T3
$ws.Range("T3").Select() | Out-Null
